$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-08-06 Tuesday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-08-07 Wednesday", 2) | Out-Null

# Update each math problem cell in the table, in document order (row-major),
# using a per-cell scoped Find/Replace so duplicate old values (e.g. "28-5=")
# are mapped to the correct new value for their specific cell.
$pairs = @(
    @("45-3=", "20+61="),
    @("50+23=", "69-31="),
    @("52-22=", "99-59="),
    @("2-0=", "12+75="),
    @("16-5=", "27+41="),
    @("31-9=", "83+5="),
    @("58+28=", "98-70="),
    @("72-0=", "50-10="),
    @("90-61=", "51+27="),
    @("78-53=", "36+3="),
    @("46+32=", "67-2="),
    @("82-60=", "34+43="),
    @("99-63=", "21+77="),
    @("13+85=", "71-71="),
    @("78-42=", "26-1="),
    @("55+22=", "68+9="),
    @("38-13=", "26+0="),
    @("80+16=", "66+1="),
    @("4+70=", "48-24="),
    @("10+60=", "50+1="),
    @("64-42=", "34+31="),
    @("65-12=", "26+42="),
    @("67+27=", "60-56="),
    @("43-13=", "53+45="),
    @("10+64=", "81-2="),
    @("86+11=", "92-59="),
    @("0+41=", "29+31="),
    @("50-26=", "37+3="),
    @("90+5=", "18+52="),
    @("98-87=", "65-56="),
    @("35+39=", "15+59="),
    @("81+16=", "97-39="),
    @("81-41=", "86-42="),
    @("63-2=", "18+48="),
    @("94-8=", "60+4="),
    @("32-18=", "41+45="),
    @("62-22=", "26+7="),
    @("37-19=", "23+43="),
    @("36+25=", "67-2="),
    @("71-27=", "30+64="),
    @("62+5=", "25-20="),
    @("68-17=", "36+15="),
    @("83+1=", "72+7="),
    @("28-5=", "91-79="),
    @("27+10=", "24+6="),
    @("29+16=", "9+16="),
    @("40-21=", "55+6="),
    @("38-25=", "45+10="),
    @("98-40=", "21-14="),
    @("58-24=", "19+61="),
    @("25-22=", "82+13="),
    @("73-18=", "95-83="),
    @("69-2=", "26+19="),
    @("60+25=", "12+77="),
    @("78-55=", "29+17="),
    @("36+56=", "55-22="),
    @("49-30=", "19+13="),
    @("28-5=", "90-59="),
    @("1+13=", "11+35="),
    @("44+40=", "32+56="),
    @("96-16=", "20+42="),
    @("23+74=", "16+68="),
    @("66+21=", "29+41="),
    @("27-19=", "14+30="),
    @("66-18=", "32+61="),
    @("69-23=", "79+14="),
    @("90-72=", "11+16="),
    @("33-32=", "87-25="),
    @("14+10=", "95-4="),
    @("85-70=", "44-37="),
    @("94+2=", "46+12="),
    @("21+58=", "89-8="),
    @("59+2=", "87-50="),
    @("30+38=", "16+54="),
    @("35+13=", "82-77="),
    @("67+28=", "61+10="),
    @("41-4=", "72-47="),
    @("16+35=", "85-17="),
    @("28-5=", "51+32="),
    @("73-12=", "26-0="),
    @("54-34=", "32-22="),
    @("15+51=", "93-56="),
    @("55-27=", "13+71="),
    @("82-36=", "68-8="),
    @("26+33=", "54-37="),
    @("37+2=", "49-0="),
    @("51-0=", "83-82="),
    @("95-33=", "7+16="),
    @("35-28=", "68-56="),
    @("83-33=", "71-47="),
    @("43-25=", "76-2="),
    @("96-4=", "73+13="),
    @("44+20=", "30+53="),
    @("45+21=", "30+6="),
    @("48+42=", "90+3="),
    @("40-8=", "5+32="),
    @("53-32=", "80-30="),
    @("68-55=", "66-23="),
    @("45+27=", "96-3="),
    @("30+34=", "22+0=")
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $old = $pairs[$idx][0]
        $new = $pairs[$idx][1]
        $cellRange = $t.Cell($r, $c).Range
        $cellRange.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
        $idx = $idx + 1
    }
}

Write-Output "Done: updated $idx cells"